$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Query" column (header + SQL-insert formulas, with their
# formatting) from E to F, opening up a blank column E for the new
# "ActiveFlg" field.
$ws.Columns("E:E").Insert()

# New header + data for column E.
$ws.Range("E1").Value = "ActiveFlg"
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 5).Formula = "1"
}

# Rebuild the formula text in column F to include the new ActiveFlg column.
for ($r = 2; $r -le 11; $r++) {
    $newFormula = '="INSERT INTO "&A' + $r + '&" ([" &B$1 &"],[" &C$1&"],[" &D$1&"],[" &E$1&"]) VALUES ( ''" & B' + $r + ' & "'',''" & C' + $r + ' & "'',''" & D' + $r + ' & "'',''" & E' + $r + ' & "'')"'
    $ws.Cells.Item($r, 6).Formula = $newFormula
}

$ws.Range("F2:F11").Select()
